$excel = New-Object -ComObject Excel.Application
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell value updates ---
$ws.Range("A2").Value = 2993253
$ws.Range("D2").Value = 216.5
$ws.Range("D3").Value = 216.5
$ws.Range("D45").Value = 98
$ws.Range("D46").Value = 228

# --- Conditional-formatting rules were trialled (via the "New Formatting
# Rule" dialog, picking theme fill swatches "Background 2, Darker 10%" and
# "White, Background 1, Darker 5%" from the color gallery, duplicating one
# rule along the way) and then removed again. Excel keeps the differential
# styles (dxfs) used by those rules around even after the rules themselves
# are deleted, so three orphaned dxf entries remain behind in the style
# table. Reproduce that by adding the conditional formats with their fills,
# then deleting the rules. ---

$rng = $ws.Range("A2:D3")

$fc1 = $rng.FormatConditions.Add(2, 0, "=TRUE")
$fc1.Interior.Color = 13553360
$rng.FormatConditions.Delete()

$fc2 = $rng.FormatConditions.Add(2, 0, "=TRUE")
$fc2.Font.Strikethrough = $false
$fc2.Interior.Color = 13553360
$rng.FormatConditions.Delete()

$fc3 = $rng.FormatConditions.Add(2, 0, "=TRUE")
$fc3.Interior.Color = 15921906
$rng.FormatConditions.Delete()

# --- View state: scroll the window and move the selection ---
$excel.ActiveWindow.ScrollRow = 37
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A2:D3").Select()

Write-Output "done"
